$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: add P1=14, Q1=15 with the same style as the existing header cells (O1) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the format from the existing header cell (O1, style index 1: bold,
# bordered, centered) onto the two new header cells instead of rebuilding
# the format property-by-property (avoids minting redundant style entries).
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 2-25: swap I/K and M/O columns (1<->2), and add new P, Q columns = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2 (new column)
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2 (new column)
}
